$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.3997546666666666
$ws.Range("H2").Value = 1.199264
$ws.Range("I2").Value = 0.1320462084214824
$ws.Range("J2").Value = 0.1320462084214824
$ws.Range("M2").Value = 7.106976666666665
$ws.Range("N2").Value = 21.32093
$ws.Range("O2").Value = 0.1598176868560746
$ws.Range("P2").Value = 0.1598176868560746
$ws.Range("Q2").Value = 2.841047088391111
$ws.Range("R2").Value = 25.56942379552
$ws.Range("S2").Value = 0.02110331958803643
$ws.Range("T2").Value = 0.02110331958803643
$ws.Range("G3").Value = 0.3997546666666666
$ws.Range("H3").Value = 1.199264
$ws.Range("I3").Value = 0.1320462084214824
$ws.Range("J3").Value = 0.1320462084214824
$ws.Range("O3").Value = 0.6350325402576649
$ws.Range("P3").Value = 0.6350325402576648
$ws.Range("Q3").Value = 11.28884659153778
$ws.Range("R3").Value = 101.59961932384
$ws.Range("S3").Value = 0.08385363916528703
$ws.Range("T3").Value = 0.08385363916528701
$ws.Range("G4").Value = 0.3997546666666666
$ws.Range("H4").Value = 1.199264
$ws.Range("I4").Value = 0.1320462084214824
$ws.Range("J4").Value = 0.1320462084214824
$ws.Range("O4").Value = 0.2051497728862606
$ws.Range("P4").Value = 0.2051497728862606
$ws.Range("Q4").Value = 3.646906524604444
$ws.Range("R4").Value = 32.82215872144
$ws.Range("S4").Value = 0.02708924966815895
$ws.Range("T4").Value = 0.02708924966815895
$ws.Range("I5").Value = 0.6840925621829359
$ws.Range("J5").Value = 0.684092562182936
$ws.Range("M5").Value = 7.106976666666665
$ws.Range("N5").Value = 21.32093
$ws.Range("O5").Value = 0.1598176868560746
$ws.Range("P5").Value = 0.1598176868560746
$ws.Range("Q5").Value = 14.71862922240222
$ws.Range("R5").Value = 132.46766300162
$ws.Range("S5").Value = 0.1093300908835222
$ws.Range("T5").Value = 0.1093300908835222
$ws.Range("I6").Value = 0.6840925621829359
$ws.Range("J6").Value = 0.684092562182936
$ws.Range("O6").Value = 0.6350325402576649
$ws.Range("P6").Value = 0.6350325402576648
$ws.Range("S6").Value = 0.4344210375344044
$ws.Range("T6").Value = 0.4344210375344044
$ws.Range("I7").Value = 0.6840925621829359
$ws.Range("J7").Value = 0.684092562182936
$ws.Range("O7").Value = 0.2051497728862606
$ws.Range("P7").Value = 0.2051497728862606
$ws.Range("S7").Value = 0.1403414337650094
$ws.Range("T7").Value = 0.1403414337650094
$ws.Range("G8").Value = 0.5566186666666667
$ws.Range("I8").Value = 0.1838612293955817
$ws.Range("J8").Value = 0.1838612293955817
$ws.Range("M8").Value = 7.106976666666665
$ws.Range("N8").Value = 21.32093
$ws.Range("O8").Value = 0.1598176868560746
$ws.Range("P8").Value = 0.1598176868560746
$ws.Range("Q8").Value = 3.955875876231111
$ws.Range("R8").Value = 35.60288288608
$ws.Range("S8").Value = 0.02938427638451598
$ws.Range("T8").Value = 0.02938427638451597
$ws.Range("G9").Value = 0.5566186666666667
$ws.Range("I9").Value = 0.1838612293955817
$ws.Range("J9").Value = 0.1838612293955817
$ws.Range("O9").Value = 0.6350325402576649
$ws.Range("P9").Value = 0.6350325402576648
$ws.Range("S9").Value = 0.1167578635579735
$ws.Range("T9").Value = 0.1167578635579735
$ws.Range("G10").Value = 0.5566186666666667
$ws.Range("I10").Value = 0.1838612293955817
$ws.Range("J10").Value = 0.1838612293955817
$ws.Range("O10").Value = 0.2051497728862606
$ws.Range("P10").Value = 0.2051497728862606
$ws.Range("Q10").Value = 5.077955097084445
$ws.Range("R10").Value = 45.70159587376
$ws.Range("S10").Value = 0.03771908945309226
$ws.Range("T10").Value = 0.03771908945309225
